$d = $word.ActiveDocument

function Find-ParagraphByText($doc, $needle) {
    $n = $doc.Paragraphs.Count
    for ($i = 1; $i -le $n; $i++) {
        $p = $doc.Paragraphs($i)
        if ($p.Range.Text -like $needle) {
            return $p
        }
    }
    return $null
}

# --- Change 1: "The SongNotes field" -> "The Notes field" (also drops the now-stale
#     spell-check proofErr wrap that surrounded the old "SongNotes" word) ---
$d.Content.Find.Execute("The SongNotes field is optional", $true, $false, $false, $false, $false, $true, 1, $false, "The Notes field is optional", 2) | Out-Null

# --- Change 2: insert "on this lane, " into the "digit in each lane" sentence ---
$d.Content.Find.Execute("should appear at this point in the song", $true, $false, $false, $false, $false, $true, 1, $false, "should appear on this lane, at this point in the song", 2) | Out-Null

# --- Change 3a: bold the sentence "Every Hold note must be followed by a Release note"
#     inside the "Hold Release" paragraph ---
$pHold = Find-ParagraphByText $d "*Every Hold note must be followed*"
$full = $pHold.Range.Text
$boldTarget = "Every Hold note must be followed by a Release note"
$idx = $full.IndexOf($boldTarget)
$bStart = $pHold.Range.Start + $idx
$bEnd = $bStart + $boldTarget.Length
$boldRange = $d.Range($bStart, $bEnd)
$boldRange.Bold = 1

# --- Change 3b: insert a new blank paragraph right before the "Note: Multiple notes..."
#     paragraph (it keeps its own text/formatting untouched, just moves down one) ---
$pNote = Find-ParagraphByText $d "Note: Multiple notes can appear*"
$pNote.Range.InsertParagraphBefore() | Out-Null

# --- Change 3c: append the new "How To Play" sentence to the end of the
#     "Typically, the first lane..." paragraph ---
$pTyp = Find-ParagraphByText $d "Typically, the first lane*"
$full2 = $pTyp.Range.Text
$tail = "notes."
$tIdx = $full2.IndexOf($tail)
$tStart = $pTyp.Range.Start + $tIdx
$tEnd = $tStart + $tail.Length
$tailRange = $d.Range($tStart, $tEnd)
$tailRange.Text = "notes. For more information on which notes should be used for each difficulty, consult the How To Play section inside the game. "

Write-Output "All changes applied"
